$d = $word.ActiveDocument

$replacements = @(
    @{old = "429×8="; new = "837×8="},
    @{old = "428×3="; new = "419×8="},
    @{old = "255×5="; new = "248×4="},
    @{old = "480×7="; new = "676×8="},
    @{old = "635×8="; new = "365×5="},
    @{old = "429×7="; new = "207×2="},
    @{old = "558×9="; new = "115×4="},
    @{old = "651×6="; new = "389×5="},
    @{old = "958×2="; new = "780×9="},
    @{old = "860×8="; new = "457×4="},
    @{old = "698×4="; new = "128×2="},
    @{old = "566×2="; new = "492×2="},
    @{old = "234×8="; new = "616×4="},
    @{old = "624×3="; new = "177×8="},
    @{old = "857×5="; new = "376×5="},
    @{old = "728×8="; new = "665×5="},
    @{old = "508×7="; new = "210×6="},
    @{old = "639×5="; new = "919×8="},
    @{old = "733×3="; new = "570×8="},
    @{old = "505×4="; new = "674×4="},
    @{old = "108×2="; new = "841×7="},
    @{old = "577×6="; new = "434×5="},
    @{old = "900×5="; new = "668×3="},
    @{old = "717×2="; new = "820×6="},
    @{old = "238×7="; new = "746×8="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
